$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 12.91660428699544
$ws.Range("C2").Value = 9.616184804558625
$ws.Range("D2").Value = 5.957635169912753
$ws.Range("E2").Value = 16.59543228894451
$ws.Range("G2").Value = 22.07251136895522
$ws.Range("H2").Value = 12.0893259132711
$ws.Range("N2").Value = 15.65497906061422
$ws.Range("O2").Value = 17.57040133161088
$ws.Range("B3").Value = 12.24639895983978
$ws.Range("C3").Value = 9.155296927913563
$ws.Range("D3").Value = 5.833720889532115
$ws.Range("E3").Value = 15.6471175078519
$ws.Range("G3").Value = 21.77778417982498
$ws.Range("H3").Value = 12.11724971743237
$ws.Range("N3").Value = 15.69889818209061
$ws.Range("O3").Value = 17.5588476430608
$ws.Range("B4").Value = 11.81657305076334
$ws.Range("C4").Value = 8.858538919031741
$ws.Range("D4").Value = 5.758100068346023
$ws.Range("E4").Value = 15.03950291501773
$ws.Range("G4").Value = 21.6061507787689
$ws.Range("H4").Value = 12.13738669748682
$ws.Range("N4").Value = 15.72768387281925
$ws.Range("O4").Value = 17.55786416055057
$ws.Range("B5").Value = 11.63700034879369
$ws.Range("C5").Value = 8.734236639339667
$ws.Range("D5").Value = 5.727448897353611
$ws.Range("E5").Value = 14.78578878234094
$ws.Range("G5").Value = 21.53865402907826
$ws.Range("H5").Value = 12.14634229889647
$ws.Range("N5").Value = 15.73987246918448
$ws.Range("O5").Value = 17.55899686375376
$ws.Range("B6").Value = 11.60692202952868
$ws.Range("C6").Value = 8.713395745799456
$ws.Range("D6").Value = 5.722370684298506
$ws.Range("E6").Value = 14.74329979673544
$ws.Range("G6").Value = 21.52759668388681
$ws.Range("H6").Value = 12.14787456369872
$ws.Range("N6").Value = 15.74192407244747
$ws.Range("O6").Value = 17.55927744775456
$ws.Range("B7").Value = 11.81416887000021
$ws.Range("C7").Value = 8.856876050741132
$ws.Range("D7").Value = 5.75768596359329
$ws.Range("E7").Value = 15.03610556947122
$ws.Range("G7").Value = 21.6052304678524
$ws.Range("H7").Value = 12.13750444436736
$ws.Range("N7").Value = 15.7278463962986
$ws.Range("O7").Value = 17.55787323281703
$ws.Range("B8").Value = 12.68942348744916
$ws.Range("C8").Value = 9.460179092093639
$ws.Range("D8").Value = 5.914842280518664
$ws.Range("E8").Value = 16.27385833166596
$ws.Range("G8").Value = 21.96901455275609
$ws.Range("H8").Value = 12.0983315773404
$ws.Range("N8").Value = 15.66974524965163
$ws.Range("O8").Value = 17.56514798953852
$ws.Range("B9").Value = 14.25365946851758
$ws.Range("C9").Value = 10.53070983381396
$ws.Range("D9").Value = 6.224615080605407
$ws.Range("E9").Value = 18.60905250733134
$ws.Range("G9").Value = 22.75173893348916
$ws.Range("H9").Value = 12.04536293806524
$ws.Range("N9").Value = 15.57021053716696
$ws.Range("O9").Value = 17.62796580999826
$ws.Range("B10").Value = 15.302916438227
$ws.Range("C10").Value = 11.24512636666486
$ws.Range("D10").Value = 6.450480517506954
$ws.Range("E10").Value = 20.26560464001165
$ws.Range("G10").Value = 23.36244207070325
$ws.Range("H10").Value = 12.02112836065546
$ws.Range("N10").Value = 15.50581772203865
$ws.Range("O10").Value = 17.70370867372096
$ws.Range("B11").Value = 15.75744789650991
$ws.Range("C11").Value = 11.55395088612533
$ws.Range("D11").Value = 6.552324540392704
$ws.Range("E11").Value = 20.97676976288928
$ws.Range("G11").Value = 23.64652943385287
$ws.Range("H11").Value = 12.01331767534353
$ws.Range("N11").Value = 15.47841163643403
$ws.Range("O11").Value = 17.74455629059249
$ws.Range("B12").Value = 15.92622269376973
$ws.Range("C12").Value = 11.66853677498055
$ws.Range("D12").Value = 6.590718930414941
$ws.Range("E12").Value = 21.24000830874481
$ws.Range("G12").Value = 23.75488372895691
$ws.Range("H12").Value = 12.01082387871279
$ws.Range("N12").Value = 15.46830431053116
$ws.Range("O12").Value = 17.76093775670293
$ws.Range("B13").Value = 15.890023828836
$ws.Range("C13").Value = 11.6439640570538
$ws.Range("D13").Value = 6.582458288872861
$ws.Range("E13").Value = 21.18358434695082
$ws.Range("G13").Value = 23.73151514967257
$ws.Range("H13").Value = 12.01134030231181
$ws.Range("N13").Value = 15.47046906966632
$ws.Range("O13").Value = 17.75736920191124
$ws.Range("B14").Value = 15.77140052187979
$ws.Range("C14").Value = 11.56342539173903
$ws.Range("D14").Value = 6.555486942464023
$ws.Range("E14").Value = 20.99854789873234
$ws.Range("G14").Value = 23.65542896132892
$ws.Range("H14").Value = 12.01310320163421
$ws.Range("N14").Value = 15.47757467700884
$ws.Range("O14").Value = 17.74588573734982
$ws.Range("B15").Value = 15.69830266755593
$ws.Range("C15").Value = 11.51378498619818
$ws.Range("D15").Value = 6.538942650468119
$ws.Range("E15").Value = 20.88441899428834
$ws.Range("G15").Value = 23.6089213842967
$ws.Range("H15").Value = 12.01424349667997
$ws.Range("N15").Value = 15.48196231321966
$ws.Range("O15").Value = 17.73897053139673
$ws.Range("B16").Value = 15.27274697141921
$ws.Range("C16").Value = 11.2246155274101
$ws.Range("D16").Value = 6.443802961404252
$ws.Range("E16").Value = 20.21827787215669
$ws.Range("G16").Value = 23.34399180394583
$ws.Range("H16").Value = 12.02170363602689
$ws.Range("N16").Value = 15.50764669123872
$ws.Range("O16").Value = 17.70116734858539
$ws.Range("B17").Value = 15.00579278440667
$ws.Range("C17").Value = 11.04305183800684
$ws.Range("D17").Value = 6.385175937696992
$ws.Range("E17").Value = 19.79878428733951
$ws.Range("G17").Value = 23.1829837384029
$ws.Range("H17").Value = 12.02710466848204
$ws.Range("N17").Value = 15.5238860351064
$ws.Range("O17").Value = 17.67960981981205
$ws.Range("B18").Value = 14.850106773271
$ws.Range("C18").Value = 10.93710104264379
$ws.Range("D18").Value = 6.351372551546965
$ws.Range("E18").Value = 19.55350893769889
$ws.Range("G18").Value = 23.09097437294959
$ws.Range("H18").Value = 12.03051359896979
$ws.Range("N18").Value = 15.533404074447
$ws.Range("O18").Value = 17.66781257529947
$ws.Range("B19").Value = 14.79702882930497
$ws.Range("C19").Value = 10.90096800106682
$ws.Range("D19").Value = 6.339914421985217
$ws.Range("E19").Value = 19.46977594881948
$ws.Range("G19").Value = 23.0599281678595
$ws.Range("H19").Value = 12.03171968275219
$ws.Range("N19").Value = 15.53665723787499
$ws.Range("O19").Value = 17.66392178773061
$ws.Range("B20").Value = 15.03443265737676
$ws.Range("C20").Value = 11.06253720572674
$ws.Range("D20").Value = 6.391425751744843
$ws.Range("E20").Value = 19.84385318472134
$ws.Range("G20").Value = 23.20006232196189
$ws.Range("H20").Value = 12.02649840825618
$ws.Range("N20").Value = 15.52213895227338
$ws.Range("O20").Value = 17.68184237960993
$ws.Range("B21").Value = 15.80633437582047
$ws.Range("C21").Value = 11.58714581957586
$ws.Range("D21").Value = 6.563414062458563
$ws.Range("E21").Value = 21.05306190051721
$ws.Range("G21").Value = 23.67775724232886
$ws.Range("H21").Value = 12.012572790936
$ws.Range("N21").Value = 15.47548024342784
$ws.Range("O21").Value = 17.74923397488892
$ws.Range("B22").Value = 16.29128807805013
$ws.Range("C22").Value = 11.91624504324614
$ws.Range("D22").Value = 6.674799765854648
$ws.Range("E22").Value = 21.80801825291463
$ws.Range("G22").Value = 23.99442735728536
$ws.Range("H22").Value = 12.00617646312298
$ws.Range("N22").Value = 15.44656407714858
$ws.Range("O22").Value = 17.79859823651951
$ws.Range("B23").Value = 16.03426582567513
$ws.Range("C23").Value = 11.7418675246386
$ws.Range("D23").Value = 6.615457492374204
$ws.Range("E23").Value = 21.4083054677985
$ws.Range("G23").Value = 23.82504698356261
$ws.Range("H23").Value = 12.00934227292949
$ws.Range("N23").Value = 15.46185295675444
$ws.Range("O23").Value = 17.77176714497698
$ws.Range("B24").Value = 15.02149145049572
$ws.Range("C24").Value = 11.05373275055239
$ws.Range("D24").Value = 6.388600512658221
$ws.Range("E24").Value = 19.82349031918519
$ws.Range("G24").Value = 23.19233935174436
$ws.Range("H24").Value = 12.02677155223966
$ws.Range("N24").Value = 15.522928241971
$ws.Range("O24").Value = 17.68083118084764
$ws.Range("B25").Value = 13.84761771225597
$ws.Range("C25").Value = 10.25355573605997
$ws.Range("D25").Value = 6.14092798251398
$ws.Range("E25").Value = 17.96117850399246
$ws.Range("G25").Value = 22.53323823922042
$ws.Range("H25").Value = 12.05712362738627
$ws.Range("N25").Value = 15.59560020016266
$ws.Range("O25").Value = 17.74588573734982
